$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 235, pushing the existing data (rows 235-314) down
# to 236-315. This mirrors what used to be a new weekly price observation
# prepended to the series for this market/product.
$ws.Rows.Item(235).Insert()

# Row 236 (the old row 235, now shifted down) already carries all of the
# categorical columns (A,B,C,E,F,G,H,I,N,O,Q,R) that the new row needs, so
# copy that row into the freshly inserted row 235 first...
$ws.Range("A236:R236").Copy()
$ws.Range("A235:R235").PasteSpecial()

# ...then overwrite the columns that actually hold new data for this entry:
# Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Precio $/Kg.
$ws.Range("D235").Value = 44627
$ws.Range("J235").Value = 130
$ws.Range("K235").Value = 5000
$ws.Range("L235").Value = 5500
$ws.Range("M235").Value = 5231
$ws.Range("P235").Value = 1744
